# Commit: "case with 380 kV done"
# Updates the pl_mw.xlsx line-flow result columns C,D,E,F,G,I,L,O for rows 2-25
# on the active worksheet (Sheet1) with the newly computed 380 kV case values.
# All other cells (A,B,H,J,K,M,N and rows 1, 26+) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns C:G, rows 2-25 ---
$data = @(
    @(0.1820850610635887, 0.02139114999639347, 0.1309150718809811, 1.211031564607481, 0.002419494034503993),
    @(0.180673226964899, 0.02150530346930246, 0.1282938749518223, 1.161875341057907, 0.002424323654764397),
    @(0.1798937303399271, 0.02158405416770925, 0.1267606846778548, 1.132458981401868, 0.002427448008444115),
    @(0.1795980346741715, 0.02161831890332344, 0.1261550105877873, 1.120663222324211, 0.002428761312758298),
    @(0.1795502602553398, 0.02162413966536114, 0.1260555916349624, 1.118716094122874, 0.002428981812263252),
    @(0.1798896536075461, 0.02158450748144958, 0.1267524390293495, 1.132299124909721, 0.00242746555774476),
    @(0.181580115414377, 0.02142870998961932, 0.1299954251432887, 1.193922996505123, 0.002421126379851074),
    @(0.185589518149186, 0.02119212000845394, 0.1369628760972148, 1.320892216726307, 0.00240995025061197),
    @(0.1889607340264092, 0.02106063777035061, 0.1424575016859748, 1.417989857857151, 0.002402495546440534),
    @(0.1905873037926256, 0.02101009120178432, 0.145039856133792, 1.463008112900752, 0.002399266604373111),
    @(0.191216647170819, 0.02099228925564844, 0.1460297154527836, 1.480178544764215, 0.002398067076580013),
    @(0.1910805105783737, 0.02099606356555483, 0.145815997728306, 1.476475097936799, 0.002398324386342368),
    @(0.1906388116123452, 0.02100859975196556, 0.1451210520098201, 1.464418263595377, 0.002399167454178554),
    @(0.1903700036578471, 0.02101645309768685, 0.1446969395864031, 1.457049153434042, 0.002399686875520409),
    @(0.1888563071407958, 0.0210641281084385, 0.1422904111266092, 1.415064963122177, 0.002402709818985439),
    @(0.1879515372131237, 0.0210957528753859, 0.1408353445563861, 1.389526968570152, 0.002404605759710912),
    @(0.1874398882634694, 0.02111481429449569, 0.1400062201813341, 1.374917948626319, 0.002405711533029857),
    @(0.1872681549185131, 0.02112141767928932, 0.1397268290449531, 1.369985252527783, 0.002406088556918402),
    @(0.1880469456635439, 0.02109229610165997, 0.1409894320030389, 1.392237269412476, 0.002404402353339313),
    @(0.1907681855156795, 0.02100488117750032, 0.145324848792562, 1.467956301614038, 0.002398919196148122),
    @(0.1926247717818228, 0.02095555810079119, 0.1482281421277918, 1.51816047640034, 0.002395470818923324),
    @(0.1916267220930337, 0.02098116610996215, 0.1466721859507629, 1.491299570416544, 0.002397298955267134),
    @(0.1880037849817739, 0.02109385616890691, 0.1409197459848599, 1.391011714737971, 0.002404494264359842),
    @(0.1844303053042466, 0.02124871243885806, 0.1350123955113567, 1.285879800927859, 0.002412840234471541)
)
$nrows = $data.Count
$ncols = $data[0].Count
$block = New-Object "object[,]" $nrows,$ncols
for ($r = 0; $r -lt $nrows; $r++) {
    for ($c = 0; $c -lt $ncols; $c++) {
        $block[$r,$c] = $data[$r][$c]
    }
}
$ws.Range("C2:G25").Value = $block

# --- Columns I:I, rows 2-25 ---
$data = @(
    @(1.163511203300047),
    @(1.124087299103635),
    @(1.100605263042212),
    @(1.091217549903973),
    @(1.089669662163644),
    @(1.100477923543778),
    @(1.149767023358578),
    @(1.252211094437115),
    @(1.331071227906207),
    @(1.367742183416567),
    @(1.381744171662561),
    @(1.378723442314637),
    @(1.368891815708295),
    @(1.362884727316185),
    @(1.328690808459157),
    @(1.307918716518444),
    @(1.296046131035439),
    @(1.292039127776491),
    @(1.310122175544649),
    @(1.371776461281314),
    @(1.412744699467552),
    @(1.390817253235127),
    @(1.30912577508937),
    @(1.223870829721321)
)
$nrows = $data.Count
$ncols = $data[0].Count
$block = New-Object "object[,]" $nrows,$ncols
for ($r = 0; $r -lt $nrows; $r++) {
    for ($c = 0; $c -lt $ncols; $c++) {
        $block[$r,$c] = $data[$r][$c]
    }
}
$ws.Range("I2:I25").Value = $block

# --- Columns L:L, rows 2-25 ---
$data = @(
    @(0.1758153550338051),
    @(0.1706748245311545),
    @(0.1676285239037725),
    @(0.1664146691010515),
    @(0.1662147690242222),
    @(0.1676120420978435),
    @(0.1740199911306775),
    @(0.1874653207916026),
    @(0.1978900075683612),
    @(0.2027534419481896),
    @(0.2046126905016905),
    @(0.2042114842439702),
    @(0.2029060504531088),
    @(0.2021087272073316),
    @(0.1975746216706398),
    @(0.1948242369595619),
    @(0.1932536868447698),
    @(0.1927238788907744),
    @(0.195115839373301),
    @(0.203289009956606),
    @(0.2087331380030264),
    @(0.2058180784955113),
    @(0.1949839726599549),
    @(0.1837328012529866)
)
$nrows = $data.Count
$ncols = $data[0].Count
$block = New-Object "object[,]" $nrows,$ncols
for ($r = 0; $r -lt $nrows; $r++) {
    for ($c = 0; $c -lt $ncols; $c++) {
        $block[$r,$c] = $data[$r][$c]
    }
}
$ws.Range("L2:L25").Value = $block

# --- Columns O:O, rows 2-25 ---
$data = @(
    @(4.26680577857752),
    @(4.106082441381432),
    @(4.010087063132971),
    @(3.971641334929473),
    @(3.96529801887624),
    @(4.009565849089824),
    @(4.210828249189092),
    @(4.627008972195654),
    @(4.946156478208479),
    @(5.094311143175787),
    @(5.150845182168837),
    @(5.13865033508938),
    @(5.098953572300843),
    @(5.074694405999878),
    @(4.936534349567978),
    @(4.852541587904398),
    @(4.804510775138681),
    @(4.788296280678253),
    @(4.861453797260026),
    @(5.110601745461167),
    @(5.275948717081974),
    @(5.187468660132424),
    @(4.857423784822913),
    @(4.512091093311255)
)
$nrows = $data.Count
$ncols = $data[0].Count
$block = New-Object "object[,]" $nrows,$ncols
for ($r = 0; $r -lt $nrows; $r++) {
    for ($c = 0; $c -lt $ncols; $c++) {
        $block[$r,$c] = $data[$r][$c]
    }
}
$ws.Range("O2:O25").Value = $block

Write-Host "Updated columns C,D,E,F,G,I,L,O for rows 2-25 on" $ws.Name
